# Update the Kilosort summary workbook:
#  - fill in newly-finished pipeline steps (Started/Finished/Discussed/etc.)
#    for sessions 210409 (S17), 210511+210512 (S20), 210912+210913 (S25)
#  - add unit counts and notes for those rows
#  - re-colour the now-complete rows to match the "all steps completed" style
#  - unhide the pipeline-step columns (E:J) that were only used for
#    intermediate tracking
#  - move the active selection to B9:D9 and drop the old frozen scroll
#    position

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Unhide columns E:J (pipeline tracking columns)
# ---------------------------------------------------------------------
$ws.Range("E1:J1").EntireColumn.Hidden = $false

# ---------------------------------------------------------------------
# 2) Row 7 -- ID 17, date 210409 -- now fully complete
# ---------------------------------------------------------------------
$ws.Range("F7:J7").Value = 1
$ws.Range("L7").Value = 26
$ws.Range("M7").Value = "a lot of the clusters were noise/purely contamination"

# Row 9 -- ID 17, date 210410 -- "Finished" step done
$ws.Range("F9").Value = 1

# ---------------------------------------------------------------------
# 3) Row 23 -- ID 20, date 210511 -- now fully complete
# ---------------------------------------------------------------------
$ws.Range("I23:J23").Value = 1
$ws.Range("L23").Value = 33

# Row 25 -- ID 20, date 210512 -- now fully complete
$ws.Range("I25:J25").Value = 1
$ws.Range("L25").Value = 49
$ws.Range("M25").Value = "3 clusters were exluded during post-curation"

# ---------------------------------------------------------------------
# 4) Row 53 -- ID 25, date 210912 -- now fully complete
# ---------------------------------------------------------------------
$ws.Range("F53:J53").Value = 1
$ws.Range("L53").Value = 42
$ws.Range("M53").Value = "maybe a little more unstable than usual? "

# ---------------------------------------------------------------------
# 5) Re-colour completed rows to match the "all steps completed" style
#    (copy full row formatting from row 56, which already uses it, and
#    the date-column accent style from C56)
# ---------------------------------------------------------------------
$ws.Range("A56:M56").Copy()
$ws.Range("A7:M7").PasteSpecial(-4122)
$ws.Range("A23:M23").PasteSpecial(-4122)
$ws.Range("A25:M25").PasteSpecial(-4122)
$ws.Range("A53:M53").PasteSpecial(-4122)
$ws.Range("A55:M55").PasteSpecial(-4122)

$ws.Range("C56").Copy()
$ws.Range("C53").PasteSpecial(-4122)
$ws.Range("C55").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 6) Move selection / scroll position
# ---------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B9:D9").Select()
